# Auto-generated edit applying the crypto-price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "98.209.60"
$ws.Range("E2").Value2 = "  +4.97%  "
$ws.Range("D3").Value2 = "3.356.98"
$ws.Range("E3").Value2 = "  +10.04%  "
$ws.Range("E4").Value2 = "  -0.13%  "
$ws.Range("D5").Value2 = "'258.09"
$ws.Range("E5").Value2 = "  +10.48%  "
$ws.Range("D6").Value2 = "'621.58"
$ws.Range("E6").Value2 = "  +2.64%  "
$ws.Range("D7").Value2 = "'1.21"
$ws.Range("E7").Value2 = "  +10.34%  "
$ws.Range("E8").Value2 = "  +2.56%  "
$ws.Range("D9").Value2 = "'0.999"
$ws.Range("E9").Value2 = "  -0.04%  "
$ws.Range("D10").Value2 = "3.355.53"
$ws.Range("E10").Value2 = "  +10.08%  "
$ws.Range("D11").Value2 = "'0.800"
$ws.Range("E11").Value2 = "  +0.04%  "
$ws.Range("E12").Value2 = "  +2.52%  "
$ws.Range("D13").Value2 = "97.937.54"
$ws.Range("E13").Value2 = "  +4.85%  "
$ws.Range("D14").Value2 = "'35.83"
$ws.Range("E14").Value2 = "  +7.17%  "
$ws.Range("D15").Value2 = "'0.0000246"
$ws.Range("E15").Value2 = "  +2.98%  "
$ws.Range("D16").Value2 = "3.962.68"
$ws.Range("E16").Value2 = "  +9.33%  "
$ws.Range("D17").Value2 = "'5.50"
$ws.Range("E17").Value2 = "  +4.69%  "
$ws.Range("D18").Value2 = "3.354.40"
$ws.Range("E18").Value2 = "  +10.45%  "
$ws.Range("E19").Value2 = "  +2.86%  "
$ws.Range("D20").Value2 = "'14.99"
$ws.Range("E20").Value2 = "  +4.67%  "
$ws.Range("D21").Value2 = "'482.68"
$ws.Range("E21").Value2 = "  +10.68%  "
$ws.Range("D22").Value2 = "'5.84"
$ws.Range("E22").Value2 = "  +3.10%  "
$ws.Range("D23").Value2 = "'0.0000205"
$ws.Range("E23").Value2 = "  +8.32%  "
$ws.Range("D24").Value2 = "'9.18"
$ws.Range("E24").Value2 = "  +4.92%  "
$ws.Range("D25").Value2 = "'5.65"
$ws.Range("E25").Value2 = "  +3.41%  "
$ws.Range("D26").Value2 = "'88.26"
$ws.Range("E26").Value2 = "  +5.01%  "
$ws.Range("D27").Value2 = "'12.06"
$ws.Range("E27").Value2 = "  +2.91%  "
$ws.Range("E28").Value2 = "  +10.30%  "
$ws.Range("E29").Value2 = "  +0.09%  "
$ws.Range("B30").Value2 = "Stellar"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value2 = "'0.252"
$ws.Range("E30").Value2 = "  +1.37%  "
$ws.Range("B31").Value2 = "Cronos"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").Value2 = "'0.184"
$ws.Range("E31").Value2 = "  +4.52%  "
$ws.Range("B32").Value2 = "Binance-PegBSC-USD"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value2 = "'0.999"
$ws.Range("E32").Value2 = "  -0.11%  "
$ws.Range("B33").Value2 = "Hedera"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value2 = "'0.121"
$ws.Range("E33").Value2 = "  +0.32%  "
$ws.Range("D34").Value2 = "'9.28"
$ws.Range("E34").Value2 = "  +2.90%  "
$ws.Range("D35").Value2 = "'27.20"
$ws.Range("E35").Value2 = "  +7.71%  "
$ws.Range("B36").Value2 = "RenderToken"
$ws.Range("C36").Value2 = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D36").Value2 = "'7.39"
$ws.Range("E36").Value2 = "  -2.88%  "
$ws.Range("B37").Value2 = "Kaspa"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value2 = "'0.151"
$ws.Range("E37").Value2 = "  -2.15%  "
$ws.Range("B38").Value2 = "Bittensor"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value2 = "'514.60"
$ws.Range("E38").Value2 = "  +12.10%  "
$ws.Range("E39").Value2 = "  +3.95%  "
$ws.Range("D40").Value2 = "'24.88"
$ws.Range("E40").Value2 = "  +4.04%  "
$ws.Range("D41").Value2 = "'0.447"
$ws.Range("E41").Value2 = "  +1.98%  "
$ws.Range("D42").Value2 = "'1.26"
$ws.Range("E42").Value2 = "  +1.96%  "
$ws.Range("D43").Value2 = "'3.61"
$ws.Range("E43").Value2 = "  -2.72%  "
$ws.Range("D44").Value2 = "'3.26"
$ws.Range("E44").Value2 = "  +4.97%  "
$ws.Range("E45").Value2 = "  -0.02%  "
$ws.Range("D46").Value2 = "'0.776"
$ws.Range("E46").Value2 = "  +17.11%  "
$ws.Range("D47").Value2 = "'160.29"
$ws.Range("E47").Value2 = "  +0.15%  "
$ws.Range("D48").Value2 = "'1.91"
$ws.Range("E48").Value2 = "  +5.20%  "
$ws.Range("B49").Value2 = "ImmutableX"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").Value2 = "'1.37"
$ws.Range("E49").Value2 = "  +7.88%  "
$ws.Range("B50").Value2 = "OKB"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").Value2 = "'45.53"
$ws.Range("E50").Value2 = "  +4.22%  "
$ws.Range("B51").Value2 = "Filecoin"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").Value2 = "'4.51"
$ws.Range("E51").Value2 = "  +7.23%  "
